$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.396.82"
$ws.Range("E2").Value = "  +0.00%  "

$ws.Range("D3").Value = "3.509.84"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'591.86"
$ws.Range("E5").Value = "  +1.07%  "

$ws.Range("D6").Value = "'134.54"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +0.31%  "

$ws.Range("D9").Value = "'7.58"
$ws.Range("E9").Value = "  +5.34%  "

$ws.Range("E10").Value = "  +0.82%  "

$ws.Range("D11").Value = "'0.390"
$ws.Range("E11").Value = "  +3.81%  "

$ws.Range("D12").Value = "4.109.68"
$ws.Range("E12").Value = "  +0.52%  "

$ws.Range("E13").Value = "  +1.16%  "

$ws.Range("E14").Value = "  +0.88%  "

$ws.Range("D15").Value = "3.510.66"
$ws.Range("E15").Value = "  +0.47%  "

$ws.Range("D16").Value = "'25.84"
$ws.Range("E16").Value = "  +2.33%  "

$ws.Range("D17").Value = "64.373.62"
$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").Value = "'10.00"
$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("D19").Value = "'13.66"
$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("D20").Value = "'5.77"
$ws.Range("E20").Value = "  +2.34%  "

$ws.Range("D21").Value = "'392.35"
$ws.Range("E21").Value = "  +1.31%  "

$ws.Range("D22").Value = "'0.582"
$ws.Range("E22").Value = "  +2.87%  "

$ws.Range("D23").Value = "3.649.88"
$ws.Range("E23").Value = "  +0.51%  "

$ws.Range("D24").Value = "'74.51"
$ws.Range("E24").Value = "  +0.64%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("D26").Value = "'5.69"
$ws.Range("E26").Value = "  -0.63%  "

$ws.Range("D27").Value = "'0.0000118"
$ws.Range("E27").Value = "  +4.11%  "

$ws.Range("E28").Value = "  +0.14%  "

$ws.Range("E29").Value = "  -0.36%  "

$ws.Range("D30").Value = "'2.28"
$ws.Range("E30").Value = "  +2.25%  "

$ws.Range("D31").Value = "'8.25"
$ws.Range("E31").Value = "  +0.68%  "

$ws.Range("E32").Value = "  -4.51%  "

$ws.Range("E33").Value = "  +7.45%  "

$ws.Range("D34").Value = "3.535.77"
$ws.Range("E34").Value = "  +0.71%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").Value = "'23.38"
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").Value = "'5.36"
$ws.Range("E37").Value = "  +1.43%  "

$ws.Range("D38").Value = "'6.96"
$ws.Range("E38").Value = "  +1.90%  "

$ws.Range("E39").Value = "  +1.80%  "

$ws.Range("D40").Value = "'166.57"
$ws.Range("E40").Value = "  +2.28%  "

$ws.Range("D41").Value = "'0.0789"
$ws.Range("E41").Value = "  +1.46%  "

$ws.Range("E42").Value = "  +1.17%  "

$ws.Range("E43").Value = "  -0.14%  "

$ws.Range("E44").Value = "  +1.08%  "

$ws.Range("D45").Value = "'25.05"
$ws.Range("E45").Value = "  -2.11%  "

$ws.Range("E46").Value = "  +0.07%  "

$ws.Range("D47").Value = "'1.18"
$ws.Range("E47").Value = "  -2.74%  "

$ws.Range("D48").Value = "'6.81"
$ws.Range("E48").Value = "  +0.90%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.372.58"
$ws.Range("E49").Value = "  -3.96%  "

$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").Value = "'0.906"
$ws.Range("E50").Value = "  +0.45%  "

$ws.Range("D51").Value = "'0.0260"
$ws.Range("E51").Value = "  +0.32%  "
